# Weekly "Fruta / Hortaliza" data refresh: a new price observation is
# inserted as row 77 (pushing the existing rows 77-138 down to 78-139),
# matching the upstream source's weekly append-at-top behaviour for this
# market/category subset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 77; everything below shifts down one.
$ws.Rows.Item(77).Insert()

# Populate the new row with the latest observation.
$ws.Range("A77").Value = 4
$ws.Range("B77").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C77").Value = "Los Lagos"
$ws.Range("D77").Value2 = 44827
$ws.Range("E77").Value = 10
$ws.Range("F77").Value = 100112022
$ws.Range("G77").Value = "Arveja Verde"
$ws.Range("H77").Value = "Perfection"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 70
$ws.Range("K77").Value = 38000
$ws.Range("L77").Value = 38000
$ws.Range("M77").Value = 38000
$ws.Range("N77").Value = "$/malla 25 kilos"
$ws.Range("O77").Value = "Provincia de Huasco"
$ws.Range("P77").Value = 1520
$ws.Range("Q77").Value = 25
$ws.Range("R77").Value = "Hortaliza"
